$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table (Date | Daily Cases | Confirmed | Deaths | Recovered | Discarded | Analyze)
# currently ends at row 104 (13-Jul-2020). Three more days of data need to be appended:
# 14-Jul-2020, 15-Jul-2020 and 16-Jul-2020 -> rows 105, 106, 107.

# First, extend the existing formatting (date format in col A, centered numbers in B:G)
# from the last populated row down into the new rows, the same way Excel does when you
# continue a table.
$ws.Range("A104:G104").Copy()
$ws.Range("A105:G107").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# New daily data rows: Date (as serial numbers), Daily Cases, Confirmed, Deaths,
# Recovered, Discarded, Analyze
$newRows = @(
    @(44026, 9, 1454, 74, 1113, 888, 2),
    @(44027, 6, 1460, 74, 1127, 895, 2),
    @(44028, 8, 1468, 75, 1142, 908, 2)
)

$startRow = 105
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $newRows[$i]
    for ($col = 0; $col -lt $rowValues.Count; $col++) {
        $ws.Cells.Item($rowNum, $col + 1).Value = $rowValues[$col]
    }
}

# Move the selection to where the user would continue typing next
$ws.Range("B108").Select()
